$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$__style = $ws.Range("D2").Style
$ws.Range("D2").Value = "'49.776.05"
$ws.Range("D2").Style = $__style
$ws.Range("E2").Value = '  -0.54%  '

$__style = $ws.Range("D3").Style
$ws.Range("D3").Value = "'2.656.63"
$ws.Range("D3").Style = $__style
$ws.Range("E3").Value = '  +0.01%  '

$ws.Range("E4").Value = '  +0.04%  '

$__style = $ws.Range("D5").Style
$ws.Range("D5").Value = "'112.82"
$ws.Range("D5").Style = $__style
$ws.Range("E5").Value = '  -1.30%  '

$__style = $ws.Range("D6").Style
$ws.Range("D6").Value = "'327.93"
$ws.Range("D6").Style = $__style
$ws.Range("E6").Value = '  +0.44%  '

$__style = $ws.Range("D7").Style
$ws.Range("D7").Value = "'0.525"
$ws.Range("D7").Style = $__style
$ws.Range("E7").Value = '  -0.80%  '

$__style = $ws.Range("D8").Style
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = $__style
$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("E9").Value = '  -1.23%  '

$__style = $ws.Range("D10").Style
$ws.Range("D10").Value = "'39.80"
$ws.Range("D10").Style = $__style
$ws.Range("E10").Value = '  -3.59%  '

$__style = $ws.Range("D11").Style
$ws.Range("D11").Value = "'19.97"
$ws.Range("D11").Style = $__style
$ws.Range("E11").Value = '  -0.74%  '

$__style = $ws.Range("D12").Style
$ws.Range("D12").Value = "'0.0818"
$ws.Range("D12").Style = $__style
$ws.Range("E12").Value = '  -1.02%  '

$ws.Range("E13").Value = '  +2.21%  '

$__style = $ws.Range("D14").Style
$ws.Range("D14").Value = "'7.58"
$ws.Range("D14").Style = $__style
$ws.Range("E14").Value = '  +2.37%  '

$__style = $ws.Range("D15").Style
$ws.Range("D15").Value = "'3.073.87"
$ws.Range("D15").Style = $__style
$ws.Range("E15").Value = '  +0.12%  '

$__style = $ws.Range("D16").Style
$ws.Range("D16").Value = "'2.659.27"
$ws.Range("D16").Style = $__style
$ws.Range("E16").Value = '  -1.15%  '

$__style = $ws.Range("D17").Style
$ws.Range("D17").Value = "'0.866"
$ws.Range("D17").Style = $__style
$ws.Range("E17").Value = '  -1.23%  '

$__style = $ws.Range("D18").Style
$ws.Range("D18").Value = "'49.739.35"
$ws.Range("D18").Style = $__style
$ws.Range("E18").Value = '  -0.49%  '

$__style = $ws.Range("D19").Style
$ws.Range("D19").Value = "'13.53"
$ws.Range("D19").Style = $__style

$__style = $ws.Range("D20").Style
$ws.Range("D20").Value = "'2.93"
$ws.Range("D20").Style = $__style
$ws.Range("E20").Value = '  -0.31%  '

$__style = $ws.Range("D21").Style
$ws.Range("D21").Value = "'6.71"
$ws.Range("D21").Style = $__style
$ws.Range("E21").Value = '  -1.09%  '

$ws.Range("E22").Value = '  -0.98%  '

$__style = $ws.Range("D23").Style
$ws.Range("D23").Value = "'269.63"
$ws.Range("D23").Style = $__style
$ws.Range("E23").Value = '  -2.67%  '

$__style = $ws.Range("D24").Style
$ws.Range("D24").Value = "'69.30"
$ws.Range("D24").Style = $__style
$ws.Range("E24").Value = '  -4.54%  '

$__style = $ws.Range("D25").Style
$ws.Range("D25").Value = "'2.57"
$ws.Range("D25").Style = $__style
$ws.Range("E25").Value = '  -1.04%  '

$__style = $ws.Range("D26").Style
$ws.Range("D26").Value = "'26.23"
$ws.Range("D26").Style = $__style
$ws.Range("E26").Value = '  -2.76%  '

$__style = $ws.Range("D27").Style
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = $__style
$ws.Range("E27").Value = '  -0.05%  '

$__style = $ws.Range("D28").Style
$ws.Range("D28").Value = "'10.22"
$ws.Range("D28").Style = $__style
$ws.Range("E28").Value = '  +1.73%  '

$ws.Range("E29").Value = '  -0.48%  '

$ws.Range("E30").Value = '  -2.10%  '

$__style = $ws.Range("D31").Style
$ws.Range("D31").Value = "'34.98"
$ws.Range("D31").Style = $__style
$ws.Range("E31").Value = '  -5.34%  '

$__style = $ws.Range("D32").Style
$ws.Range("D32").Value = "'49.55"
$ws.Range("D32").Style = $__style
$ws.Range("E32").Value = '  -1.23%  '

$__style = $ws.Range("D33").Style
$ws.Range("D33").Value = "'5.50"
$ws.Range("D33").Style = $__style
$ws.Range("E33").Value = '  +0.07%  '

$ws.Range("E34").Value = '  +0.49%  '

$__style = $ws.Range("D35").Style
$ws.Range("D35").Value = "'19.16"
$ws.Range("D35").Style = $__style
$ws.Range("E35").Value = '  -1.84%  '

$ws.Range("E36").Value = '  -0.07%  '

$ws.Range("E37").Value = '  -1.17%  '

$ws.Range("E38").Value = '  -1.47%  '

$__style = $ws.Range("D39").Style
$ws.Range("D39").Value = "'3.15"
$ws.Range("D39").Style = $__style
$ws.Range("E39").Value = '  +0.67%  '

$__style = $ws.Range("D40").Style
$ws.Range("D40").Value = "'23.73"
$ws.Range("D40").Style = $__style
$ws.Range("E40").Value = '  +5.98%  '

$__style = $ws.Range("D41").Style
$ws.Range("D41").Value = "'128.10"
$ws.Range("D41").Style = $__style
$ws.Range("E41").Value = '  +2.79%  '

$ws.Range("E42").Value = '  +8.60%  '

$__style = $ws.Range("D43").Style
$ws.Range("D43").Value = "'2.31"
$ws.Range("D43").Style = $__style
$ws.Range("E43").Value = '  +4.24%  '

$ws.Range("E44").Value = '  -0.64%  '

$ws.Range("E45").Value = '  +0.39%  '

$__style = $ws.Range("D46").Style
$ws.Range("D46").Value = "'2.063.63"
$ws.Range("D46").Style = $__style
$ws.Range("E46").Value = '  -1.71%  '

$__style = $ws.Range("D47").Style
$ws.Range("D47").Value = "'2.12"
$ws.Range("D47").Style = $__style
$ws.Range("E47").Value = '  +6.60%  '

$ws.Range("E48").Value = '  -2.32%  '

$ws.Range("E49").Value = '  -1.46%  '

$ws.Range("E50").Value = '  -1.18%  '

$__style = $ws.Range("D51").Style
$ws.Range("D51").Value = "'59.35"
$ws.Range("D51").Style = $__style
$ws.Range("E51").Value = '  -1.85%  '
